$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellAddr, $val)
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" '41.996.25'
Set-TextValue "E2" '  +4.57%  '

Set-TextValue "D3" '2.231.89'
Set-TextValue "E3" '  +1.16%  '

Set-TextValue "E4" '  +0.00%  '

Set-TextValue "D5" '232.59'
Set-TextValue "E5" '  +1.71%  '

Set-TextValue "D6" '0.629'
Set-TextValue "E6" '  -0.43%  '

Set-TextValue "D7" '61.48'
Set-TextValue "E7" '  -4.21%  '

Set-TextValue "E8" '  +0.08%  '

Set-TextValue "D9" '0.406'
Set-TextValue "E9" '  +1.84%  '

Set-TextValue "E10" '  +1.03%  '

Set-TextValue "D11" '0.0906'
Set-TextValue "E11" '  +5.09%  '

Set-TextValue "E12" '  +0.07%  '

Set-TextValue "D13" '2.566.31'
Set-TextValue "E13" '  +1.24%  '

Set-TextValue "D14" '15.67'
Set-TextValue "E14" '  -2.20%  '

Set-TextValue "D15" '22.34'
Set-TextValue "E15" '  +0.49%  '

Set-TextValue "D16" '0.804'
Set-TextValue "E16" '  -1.99%  '

Set-TextValue "D17" '5.61'
Set-TextValue "E17" '  +0.10%  '

Set-TextValue "D18" '2.247.42'
Set-TextValue "E18" '  +1.47%  '

Set-TextValue "D19" '41.883.35'
Set-TextValue "E19" '  +4.39%  '

Set-TextValue "D22" '6.02'
Set-TextValue "E22" '  -1.45%  '

Set-TextValue "D23" '251.94'
Set-TextValue "E23" '  +8.06%  '

Set-TextValue "D24" '1.00'
Set-TextValue "E24" '  -0.01%  '

Set-TextValue "D25" '2.37'
Set-TextValue "E25" '  +1.80%  '

Set-TextValue "D26" '2.40'
Set-TextValue "E26" '  +1.64%  '

Set-TextValue "D27" '9.71'
Set-TextValue "E27" '  -0.31%  '

Set-TextValue "D28" '0.144'
Set-TextValue "E28" '  +1.85%  '

Set-TextValue "D29" '169.08'
Set-TextValue "E29" '  -1.75%  '

Set-TextValue "D30" '20.08'
Set-TextValue "E30" '  -0.27%  '

Set-TextValue "D31" '1.43'
Set-TextValue "E31" '  -2.28%  '

Set-TextValue "D32" '2.72'
Set-TextValue "E32" '  -0.70%  '

Set-TextValue "E33" '  -0.56%  '

Set-TextValue "D34" '5.06'
Set-TextValue "E34" '  +6.32%  '

Set-TextValue "D35" '4.68'
Set-TextValue "E35" '  +1.89%  '

Set-TextValue "D36" '0.0638'
Set-TextValue "E36" '  +1.89%  '

Set-TextValue "D37" '6.66'
Set-TextValue "E37" '  -6.00%  '

Set-TextValue "D38" '3.73'
Set-TextValue "E38" '  -4.56%  '

Set-TextValue "D39" '2.36'
Set-TextValue "E39" '  -3.61%  '

Set-TextValue "D40" '0.000260'
Set-TextValue "E40" '  +32.37%  '

Set-TextValue "E41" '  +0.03%  '

Set-TextValue "D42" '0.0241'
Set-TextValue "E42" '  +4.93%  '

Set-TextValue "D43" '4.77'
Set-TextValue "E43" '  -5.65%  '

Set-TextValue "E44" '  +2.85%  '

Set-TextValue "D45" '1.23'
Set-TextValue "E45" '  +0.31%  '

Set-TextValue "D48" '1.482.81'
Set-TextValue "E48" '  -2.70%  '

Set-TextValue "D49" '16.60'
Set-TextValue "E49" '  -5.48%  '

Set-TextValue "E50" '  +0.45%  '

Set-TextValue "D51" '52.94'
Set-TextValue "E51" '  +5.89%  '

# Row 20/21 swap: Litecoin <-> ShibaInu
Set-TextValue "B20" 'ShibaInu'
Set-TextValue "C20" 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D20" '0.0₃0908'
Set-TextValue "E20" '  -0.21%  '

Set-TextValue "B21" 'Litecoin'
Set-TextValue "C21" 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue "D21" '72.50'
Set-TextValue "E21" '  -0.10%  '

# Row 46/47 swap: Cronos <-> Aave
Set-TextValue "B46" 'Aave'
Set-TextValue "C46" 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue "D46" '99.48'
Set-TextValue "E46" '  -4.32%  '

Set-TextValue "B47" 'Cronos'
Set-TextValue "C47" 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D47" '0.0965'
Set-TextValue "E47" '  +3.82%  '

